$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Data changes (add angular-filter, add projects and admin panel) ---
# GitHub commit-day / commit-count scores bumped by one
$ws.Range("C8").Value = 17
$ws.Range("C9").Value = 41

# User Dashboard score reduced
$ws.Range("C16").Value = 5

# Admin panel: "->Projects" and "->Add Project" now scored
$ws.Range("C29").Value = 10
$ws.Range("C30").Value = 15

# --- Update the view's active selection to reflect where the author was working ---
$ws.Range("B10:E10").Select()
